# Update the dSF (column F) values for specific rows to reflect the
# "repull data, push all data, mean calculation" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -1
    "F4"  = -4
    "F10" = -8
    "F11" = -8
    "F13" = 0
    "F14" = -4
    "F15" = -3
    "F17" = 13
    "F19" = 3
    "F20" = -5
    "F23" = -6
    "F26" = -6
    "F28" = -6
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
